$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '68.219.88'
$ws.Range("E2").Value = '  -0.30%  '
$ws.Range("D3").Value = '2.705.69'
$ws.Range("E4").Value = '  +0.08%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '608.72'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.87%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '166.10'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +4.31%  '
$ws.Range("E7").Value = '  +0.05%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.553'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +2.27%  '
$ws.Range("D9").Value = '2.704.73'
$ws.Range("E9").Value = '  +2.14%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.144'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +1.18%  '
$ws.Range("E11").Value = '  +0.72%  '
$ws.Range("E12").Value = '  +2.54%  '
$ws.Range("E13").Value = '  -0.01%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '28.29'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +0.50%  '
$ws.Range("D15").Value = '3.198.68'
$ws.Range("E15").Value = '  +2.06%  '
$ws.Range("E16").Value = '  -0.44%  '
$ws.Range("D17").Value = '68.194.74'
$ws.Range("E17").Value = '  -0.23%  '
$ws.Range("D18").Value = '2.709.09'
$ws.Range("E18").Value = '  +0.78%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '11.74'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.61%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '369.11'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +1.28%  '
$ws.Range("E21").Value = '  +1.90%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.47'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +1.16%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '4.90'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +2.42%  '
$ws.Range("E24").Value = '  -2.29%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '72.78'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -3.17%  '
$ws.Range("E26").Value = '  +0.02%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.96'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -0.04%  '
$ws.Range("D28").Value = '2.843.02'
$ws.Range("E28").Value = '  +1.54%  '
$ws.Range("E29").Value = '  +0.75%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.00'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +0.63%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '577.15'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +0.43%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '8.10'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +0.31%  '
$ws.Range("E33").Value = '  +0.32%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.97'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +5.39%  '
$ws.Range("E35").Value = '  +1.06%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.999'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +0.00%  '
$ws.Range("E37").Value = '  -3.03%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '19.80'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +0.71%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '158.75'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -1.25%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.377'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +1.91%  '
$ws.Range("E41").Value = '  +0.90%  '
$ws.Range("E42").Value = '  -1.31%  '
$ws.Range("E43").Value = '  +0.80%  '
$ws.Range("E44").Value = '  -2.07%  '
$ws.Range("E45").Value = '  +0.01%  '
$ws.Range("D46").Value = '0.0₆0308'
$ws.Range("E46").Value = '  -3.72%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '40.73'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +1.12%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.595'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +3.23%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '154.75'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -2.37%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '3.88'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +1.58%  '
$ws.Range("E51").Value = '  +3.27%  '
